# Rename the picture shapes in the header/footer logos, per the commit:
#   footer (primary)      : PearsonLogo  image1.png -> image2.png
#   footer (first page)   : PearsonLogo  image1.png -> image2.png
#   header (first page)   : BTec_Logo-Orange image2.jpg -> image1.jpg
#
# InlineShape has no writable .Name in the Word object model, so each
# picture is briefly converted to a floating Shape (where .Name is
# writable), renamed, then converted back to an inline picture so the
# layout/anchoring is left exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($shapeRange, $newName) {
    if ($shapeRange.Count -gt 0) {
        $inlineShape = $shapeRange.Item(1)
        $floating = $inlineShape.ConvertToShape()
        $floating.Name = $newName
        [void]$floating.ConvertToInlineShape()
    }
}

# Primary footer (footer2.xml) - Pearson logo
Rename-InlinePicture $sec.Footers(1).Range.InlineShapes "image2.png"

# First-page footer (footer1.xml) - Pearson logo
Rename-InlinePicture $sec.Footers(2).Range.InlineShapes "image2.png"

# First-page header (header1.xml) - BTEC logo
Rename-InlinePicture $sec.Headers(2).Range.InlineShapes "image1.jpg"

Write-Output "Renamed inline picture shapes"
